$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed export has one fewer data row than before (404 -> 403),
# since two low rows were consolidated into one new top row. Remove one
# row from the data block first so the final row count lines up, then
# overwrite the whole data block (rows 2..404) with the refreshed,
# re-sorted export (descending by Saldo).
$ws.Rows.Item(3).Delete() | Out-Null

$data = @"
005268516|LUIS|250781.59
004550415|DIOGO|23638.56
004971448|CLOVIS|18977.22
002636063|LEDA|15548.77
005055226|DANILO|11717.92
004413523|ROSANE|11060.4
004570632|FABRICIO|7254.05
004240400|ADRIANA|6694.32
004547722|MARCIA|2193.01
001759765|NATAL|1985.82
004425261|THAYSA|992.79
004581652|CINCO|972.53
004813134|MONICA|961.99
004329030|DANIELA|936.81
004482102|NATALIA|936.47
004453132|BRUNO|933.72
005198093|ANA|924.01
004392159|RODRIGO|900.21
005616259|MARIA|879.74
004499920|FABIANO|877.74
004454491|MARCO|865.62
004556974|KELLY|863.5
005324981|JO|861.72
005599726|JORGE|840.7
004472538|RODOLFO|825.7
004693308|LAURA|813.54
004504449|KELMA|800
005046919|MARIANA|740
004526450|MSD|653.78
004468717|HELOISA|633.92
004462930|WALTER|608
004364200|BLOCO|574.85
005660459|MOACIR|568.58
004322719|GISELA|533.15
005293480|WAGNER|532.44
004405234|JULIO|522.09
004565108|GUSTAVO|516.91
004385806|ANILSON|504
004332747|LOHRAN|500
004460491|PEDRO|496.08
004455356|MARCELO|487.35
001761119|BLUEMETRIX|479.91
004376145|LUCYENE|466.63
004935287|ODILON|463.07
004368994|CRISTINA|453.69
005002390|LUCIANO|443.26
004508526|CASSIO|438.09
005002457|ROSANGELA|434.59
003249855|MARINA|429.81
005255637|PATRICIA|416.79
004452597|LARA|412.22
005639338|CLOTILDE|411.74
004551472|DIEGO|402.81
004479287|ANA|392.83
004515548|FLAVIA|373.31
004948033|GUILHERME|360.8
004912314|FABRICIO|357.85
004927044|CINTIA|351.77
004342617|JURACI|350.34
004862677|RENATO|338.09
005061124|BRUNO|336.38
004482163|CIRIA|336.37
004474776|GILSON|334.37
004630773|NABOR|326.03
004272426|RODRIGO|325.49
004222784|RAFAEL|317.85
004212581|MARIA|310.09
004751154|CATARINE|309.33
004212438|KENIA|300
004452912|BRUNO|300
004472404|DILSON|300
005651698|JOAO|293.99
005429528|BIANCA|292.22
004500087|JANNE|288.73
004805269|CLISIA|288.32
004754056|BRUNO|284.75
004767746|ISABELE|281.2
004565146|GUSTAVO|279.75
004487016|ROGERIO|272.99
004376853|ALBERTO|267.4
005366255|RAPHAELA|245.44
004752494|SERGIO|234.98
004332783|IRON|231.7
004165515|MAURO|229.58
004539779|RICARDO|227.8
005105970|VERA|223
004221638|CAROLINE|222.8
004584517|CAIO|222.8
004431689|LUCA|222.55
004455489|SOLANGE|217.13
004419782|NADY|214
004398253|EULER|207.76
004813088|JULIANA|201.45
002823185|SIMONE|200.4
004213929|RODOLFO|200.02
002697806|CLAUDIA|200
004395314|MARIA|200
004450760|SILVIO|200
004582648|ANUAR|200
004813166|VENIA|200
004855960|CLERIA|200
005003629|ANDRE|200
005206566|LEVI|200
005277762|NIVALDO|200
005645211|AGUINALDO|200
005683532|SYLVERSON|200
004444380|MARCELO|194.36
004213943|ELISA|176.64
004218542|JOSE|172.84
004211922|CARLOS|169.08
004508504|FERNANDO|163.81
005324840|PEDRO|161.78
004273239|DANIEL|158.26
004207184|CRISTINA|156.63
004467884|ANA|156.01
004804036|LUCIANA|154.84
005591536|GUSTAVO|153.05
005024046|ALEXANDRE|150.15
004211368|ILTON|144.16
004260002|ERICA|138.61
004332544|CELIA|129.49
004452790|GUSTAVO|129.48
004471893|PAULA|129.48
005135281|RAFAEL|128.21
004208733|REINALDO|122.81
005079458|JONAS|121.7
005103059|WALQUIRIA|120.05
004462543|RODOLFO|119.34
004895776|FERNANDO|119.2
004508516|EDUARDO|117.99
004500804|RAFAEL|117.84
005140667|MATEUS|117.47
005905737|CAIRO|115.07
005055239|NORMAN|114.08
004357848|AURELIO|109.71
000772433|MARCELO|108.32
004419765|WALTER|107.86
004536602|TATIANY|107.77
004448303|NASSIM|106
004512434|CAIO|104.34
004237325|RICARDO|103.88
002687737|JOSE|100.44
004210959|ANA|100.01
004238436|DIEGO|100.01
000330949|RENATO|100
000678704|LUDMILA|100
001000882|AYRTON|100
002973105|DARLAN|100
003641655|MARCELO|100
004202332|TATIANA|100
004207641|MAGALI|100
004211807|EDINARDO|100
004239624|NINA|100
004240014|ISABELE|100
004241147|ANTONIO|100
004252768|ALESSANDRO|100
004355790|MINEIA|100
004380749|ELEUSE|100
004404248|PAULO|100
004404342|ADSON|100
004432935|JOSE|100
004452476|IVONE|100
004458563|LUIZ|100
004458624|PEDRO|100
004461070|EDUARDO|100
004463586|MARCIA|100
004467016|ISABEL|100
004472431|LUIS|100
004473829|LUCAS|100
004479463|HENRIQUE|100
004515341|BRUNO|100
004556150|MARINA|100
004567880|LUANA|100
004588677|RACHEL|100
004756981|MATEUS|100
004805273|CLISIA|100
004827315|GILTON|100
004853111|MARCONDES|100
004866753|GENESI|100
004911541|TIAGO|100
004946542|ALESSANDRA|100
005009026|EDMUR|100
005079311|JOVINO|100
005105172|VALDIVINO|100
005135532|FELIPE|100
005304669|EMILSON|100
005313179|MAIRA|100
005379541|JOSE|100
005535788|EMILIA|100
005547702|NATHALIA|100
005554830|PAULO|100
005607508|CARLOS|100
005637820|GUILHERME|100
005698866|HELIO|100
005725431|BRUNO|100
005749972|ALESSANDRA|100
005880628|PEDRO|100
005905713|NEILA|100
005920340|PAULA|100
004563252|FERNANDO|99.42
005152037|RODRIGO|98.92
004514241|ANDRE|98.91
005044389|CLAUDIA|98.64
004946997|EDUARDO|98.52
005514036|ANA|98.27
005035754|JOSE|97
004211911|ZENILDA|96.29
004457389|RAFAEL|95.99
004975924|SERGIO|95.4
004752519|MARCUS|95.28
005312963|ALAN|94.84
004451652|MATEUS|94.75
004431591|MARIO|94.24
004383268|LAURA|93.29
004383190|MAFALDA|93.28
004267044|PATRICIA|92.61
005726697|FERNANDO|92.32
005256849|SANDRO|91.77
004472760|SANDRA|91.54
004877741|LUIZ|90.41
004479734|RODRIGO|89.08
001090818|MARIA|87.57
004381194|ALINNE|87.42
004212132|JOAO|86.38
005681354|MATHEUS|86.15
005121919|JORGE|86.09
004809902|PEDRO|85.9
004907688|HEITOR|85.27
005666419|GABRIEL|85.27
004643880|GABRIEL|84.47
005332650|DILSON|84.41
004477812|DIEGO|83.53
004332207|IRACY|83.12
004332103|JOSE|82.75
005624274|CLAYTON|82.72
004994036|BALTASAR|81.19
004318604|RENAN|80.63
004466221|WALTER|79.5
005256036|LUIZ|79.04
005032151|ANA|77.8
005558076|ALEXANDRE|75.64
005262440|BERNARDO|74.8
004321092|DANIEL|73.7
005040864|ANDRE|73.58
004340984|RENATA|72.08
004207278|CESAR|71.3
003836362|ISABELLA|70.92
004384131|ANDRE|70.85
004267119|ANA|70
004520100|ALEXANDRE|69.77
005000645|ANTONIO|68.27
004425965|CAROLLINA|65.72
004382902|LEILA|64.55
004452507|DANIELA|64.43
004216434|JAIME|64.24
001804114|WAGNER|63.98
004335251|EDMUNDO|62.82
002418050|SANDRA|62.24
004027477|GABRIELA|61.81
004691225|ANNA|61.51
004228456|FLASH|59.27
005611045|PAULA|58.6
004517506|LUIZ|58.21
005068961|JORGE|56.46
004497825|PRISCILLA|55.56
004755083|EVALDO|54.08
004552021|MARIA|53.8
004426743|GABRIELLE|53.79
004242237|MARIAH|51.94
004400640|FELIPE|51.44
005684392|PAULO|50.88
004208447|LEILA|50
004974089|CELIA|48.88
004224011|THOMAS|47.27
004374891|RODRIGO|47.04
004752615|LUZIMAR|46.2
005266369|EG|45.52
005009922|ANA|45.39
005142661|SABRINA|45.36
001731007|GUILHERME|44.73
004259659|BENTO|44.09
004265173|JULIA|43.27
005245032|ROSA|39.91
004238164|DANIELA|39.09
005685089|CARNEIRO|38.26
002401479|JULIO|37.84
005076418|LEONARDO|37.64
005924958|TIAGO|37.1
005000656|LUCIA|36.71
004419141|PAULO|34.98
005696533|CARLOS|32.82
004207374|ANGELICA|31.25
004713953|ALESSANDRA|30.85
005381719|MARIA|30.78
004398174|DANIELE|30.31
005927101|SIMONE|30
005141215|KARINA|29.6
004377415|ANGELA|28.73
005165116|ANA|28.45
003497496|ELISANDRA|28.36
004240292|MARCO|28.12
004413537|CLAUDIA|28.08
004334158|LEONE|27.65
005070742|JUSCELINO|26.98
001719494|LUIS|26.16
004224405|MILA|25.63
004350197|GISELA|25.08
004404724|LEANDRO|24.14
004886366|RENATO|22.47
005305965|SIDMAR|22.35
004371857|NAZARETH|21.52
004243043|SUELI|21.51
004388077|WLADMIR|20.89
004214604|MARIA|20.75
004224815|GUILHERME|20.14
004754920|LUIS|19.67
005143579|GABRIEL|19.2
004204255|AMADO|18.77
004556853|MARCEL|18.7
005142624|RODRIGO|17.89
004805133|PATRICIA|17.28
002894447|JOAO|16.48
004331477|SUZY|15.42
004749928|NADY|15.38
001879977|THAISSA|14.73
004422594|WANDIR|14.67
004329441|HELENA|13.74
004230529|LAIS|13
002064834|RAFAELA|12.85
004472076|RUBENS|11.8
004335031|EDMUNDO|11.17
005374916|MARCO|10.89
004459461|INTERLAGOS|10.03
004920447|MARILIA|10.03
005660155|CAROLINA|10
004421636|PATRICIA|9.88
004216298|FLORDELIZ|9.78
004264780|MARCELO|9.74
004335144|EDMUNDO|8.62
005708197|CALIDAD|8.19
004508159|FELIPE|7.63
004308815|ZELI|7.54
004752461|SERGIO|7.52
004481463|MARA|7.22
004313254|GUSTAVO|7.12
004530494|ROSANGELA|6.99
004289402|LARISSA|6.31
004482262|RITA|6.23
004448501|JOAO|5.55
004431546|GABRIELA|5.13
004346716|TIAGO|4.8
004374943|LEONARDO|4.79
004267976|E3|4.73
000834301|MARCUS|4.4
004220849|DULCE|3.59
004228090|GUSTAVO|2.98
005022526|ALEXANDRE|1.7
005142592|ALBERTO|1.58
001882235|LAGO|1.29
000431814|GUILHERME|1.1
004360430|VIOMAR|1
004486497|ELENA|0.96
004480970|ALBERTO|0.94
004115403|HEBERT|0.88
004424671|LUISA|0.88
004424761|PEDRO|0.84
004424714|HELENA|0.8
004223502|BRUNA|0.78
004587511|CARLOS|0.69
004207658|ROBERTO|0.66
004380948|LUISA|0.62
004473942|DAIANNE|0.62
004546050|LUIS|0.5
003894173|ANDREA|0.48
004806286|VERA|0.45
004453302|ISABELLA|0.39
005440756|VALERIA|0.37
004638738|GABRIEL|0.29
004278033|DAISY|0.21
005662526|AGUINALDO|0.18
004432455|LUCIANA|0.17
004612043|YURI|0.17
005075382|NAYARA|0.17
002694089|VITOR|0.16
004357159|JOAO|0.15
004320840|NATALIA|0.14
004466350|RAQUEL|0.11
004451996|ADRIANO|0.09
005047946|GABRIEL|0.09
005437764|EVA|0.09
005530256|CAROLINA|0.07
001368670|THIAGO|0.05
004223226|YESHUA|0.04
004281300|FRANKLIN|0.04
005274028|RAFAEL|0.04
004329229|GABRIEL|0.03
004213373|ALEXANDRE|0.02
004339183|JALISON|0.02
004870976|HFR|0.02
000938440|BASE|0.01
002878817|GUILHERME|0.01
004400000|VILMA|0.01
005654767|DIEGO|0.01
"@

$lines = $data -split "`n"
$startRow = 2
$endRow = $startRow + $lines.Count - 1

# Keep the "Conta" account-number column as text so leading zeros survive.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $lines.Count; $i++) {
  $parts = $lines[$i].Split('|')
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $parts[0]
  $ws.Cells.Item($r, 2).Value = $parts[1]
  $ws.Cells.Item($r, 3).Value = [double]$parts[2]
}
